$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
# "Ready for handoff" -> "Handed back: in sync with en-US" (shared string
# reused by E2/F2/E3/F3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Status column (C) shares the same "Ready for handoff" string
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
# Latest Handback DateTime refreshed (K2/K3 share the same shared string)
$wsZhCn.Range("K2").Value = "2016-09-06 05:38:11"
$wsZhCn.Range("K3").Value = "2016-09-06 05:38:11"
# Error Detail cleared now that the handback is in sync
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-09-06 05:38:27"
$wsDeDe.Range("K3").Value = "2016-09-06 05:38:27"
$wsDeDe.Range("P2").Value = ""
